# Update the "Periodo Mora" column (E16:E21) of the account statement table.
# The previous periods (2012, 2101, 2102, 2103, 2104, 2105) are replaced by the
# new set of periods, entered with the most recent period first:
# (2105, 2104, 2103, 2102, 2101, 2012).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2105"
$ws.Range("E17").Value = "2104"
$ws.Range("E18").Value = "2103"
$ws.Range("E19").Value = "2102"
$ws.Range("E20").Value = "2101"
$ws.Range("E21").Value = "2012"
